$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header row: "<name>_old" -> "<name>_FV2310" and "<name>_new" -> "<name>_FV2404"
# ---------------------------------------------------------------------------
$headers = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310",
    "diff",
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# ---------------------------------------------------------------------------
# 2. Freeze the header row
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# ---------------------------------------------------------------------------
# 3. Turn A1:U84 into a native Excel table ("Table1") with an autofilter and
#    row stripes, keeping the header row's pre-existing look (bold / shaded /
#    bordered) instead of the formatting the table style would otherwise
#    stamp on top of it.
# ---------------------------------------------------------------------------
$headerRange = $ws.Range("A1:U1")
$stashRange  = $ws.Range("A200:U200")

# Stash the header's current formatting out of the way ...
$headerRange.Copy()
$stashRange.PasteSpecial(-4122)   # xlPasteFormats

# ... clear it so adding the table doesn't need to preserve it as a one-off
# table dxf ...
$headerRange.ClearFormats()

$range = $ws.Range("A1:U84")
$listObject = $ws.ListObjects.Add(1, $range, [System.Reflection.Missing]::Value, 1)
$listObject.Name = "Table1"
$listObject.TableStyle = ""
$listObject.ShowTableStyleRowStripes = $true
$listObject.ShowTableStyleColumnStripes = $false
$listObject.ShowTableStyleFirstColumn = $false
$listObject.ShowTableStyleLastColumn = $false

# ... and finally restore the original header look in one shot.
$stashRange.Copy()
$headerRange.PasteSpecial(-4122)  # xlPasteFormats
$stashRange.Clear()

$ws.Range("A1").Select()

Write-Host "done"
